$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 20 and 21 (KB제28호스팩, 아이씨티케이), shifting remaining rows up.
$ws.Rows("20:21").Delete()
